$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: num_customers 106 -> 110, retention_rate recalculated (110/204)
$ws.Range("C21").Value = 110
$ws.Range("E21").Value = 0.5392156862745098

# Row 22: num_customers 59 -> 61, cohort_size 59 -> 61 (retention_rate stays 1)
$ws.Range("C22").Value = 61
$ws.Range("D22").Value = 61
